$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row (original row 26) entirely; remaining rows shift up by one.
$ws.Rows(26).Delete()

# Remove the "SC 92" row. After the previous delete it now sits at row 27.
$ws.Rows(27).Delete()

# Apply the remaining per-cell value changes (rows numbered as they stand
# after the two row deletions above).
$ws.Range("F6").Value = 16.43
$ws.Range("F8").ClearContents()
$ws.Range("F12").Value = 17.45
$ws.Range("F14").ClearContents()
$ws.Range("F17").Value = 17.78
$ws.Range("F18").Value = 18.35
$ws.Range("F19").ClearContents()
$ws.Range("F20").ClearContents()
$ws.Range("F23").Value = 16.48

$ws.Range("B27").Value = -20.4
$ws.Range("F27").ClearContents()
$ws.Range("B28").ClearContents()
$ws.Range("B29").ClearContents()
$ws.Range("B30").Value = -19.7
$ws.Range("B32").ClearContents()
